$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new header cells F1:H1 with the same style as the existing header row (e.g. E1)
$ws.Range("E1").Copy() | Out-Null
$ws.Range("F1:H1").PasteSpecial(-4122) | Out-Null
$ws.Application.CutCopyMode = $false

$ws.Range("F1").Value = "KNN_Outliers_MAD"
$ws.Range("G1").Value = "SVM_Outliers_MAD"
$ws.Range("H1").Value = "RF_Outliers_MAD"

# Boolean outlier flags for rows 2-8
$flags = @{
    2 = @($false, $false, $false)
    3 = @($false, $true,  $false)
    4 = @($false, $false, $false)
    5 = @($false, $true,  $false)
    6 = @($false, $false, $false)
    7 = @($false, $false, $false)
    8 = @($false, $false, $false)
}

foreach ($row in $flags.Keys) {
    $vals = $flags[$row]
    $ws.Range("F$row").Value = $vals[0]
    $ws.Range("G$row").Value = $vals[1]
    $ws.Range("H$row").Value = $vals[2]
}
